$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the delivery boy data values
$ws.Range("A1").Value = "Mark"
$ws.Range("A2").Value = "mark@gmail.com"
$ws.Range("A3").Value = 9835615595
$ws.Range("A4").Value = "Kerala"
$ws.Range("A5").Value = "merk66"
$ws.Range("A6").Value = "mark"

# Move active selection to A6 (matches final cursor position in the file)
$ws.Range("A6").Select()
